$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add H11 = "David" (shared string index 15)
$ws.Range("H11").Value = "David"

# Row 13: apply the same yellow-fill style used by A9:F9 / A10:F10 / A12:F12
# directly (avoids Copy/PasteSpecial side effects on sheet relationships).
$ws.Range("A13:F13").Interior.Color = 65535

# Add I13 = "fini" (shared string index 18)
$ws.Range("I13").Value = "fini"

# Add I14 = "fini" (shared string index 18)
$ws.Range("I14").Value = "fini"

# Update the view: active selection moves to H11 (also resets the scrolled
# topLeftCell back to the top-left of the sheet).
$ws.Range("H11").Select() | Out-Null
